$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.784.12"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.656.10"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.08"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3814"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3619"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.34"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.257"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08222"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.546"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.465"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001243"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.651.40"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06982"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.782"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.78"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("B24").Value = "WrappedBTC"
$ws.Range("C24").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.769.07"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.568"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.077"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.31"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.76"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.226"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.49"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.824.19"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.952"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.184"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.084"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.92"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.11%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.159"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2520"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08825"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.07164"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.23"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +9.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7079"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.344"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.01"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6556"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.334"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9998"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.957"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07975"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.64"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.195"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.06%  "
